$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 325 - this shifts the existing
# row 325 (and everything below it, through the former row 452) down
# by one, becoming rows 326..453, and preserves all formatting
# (including the date-number-format style on column D).
$ws.Rows.Item(325).Insert()

# Populate the newly inserted row 325 with the new record's data.
$ws.Range("A325").Value = 10
$ws.Range("B325").Value = "Vega Modelo de Temuco"
$ws.Range("C325").Value = "La Araucanía"
$ws.Range("D325").Value = 44809
$ws.Range("E325").Value = 9
$ws.Range("F325").Value = 100114014
$ws.Range("G325").Value = "Betarraga"
$ws.Range("H325").Value = "Sin especificar"
$ws.Range("I325").Value = "Primera"
$ws.Range("J325").Value = 30
$ws.Range("K325").Value = 12000
$ws.Range("L325").Value = 12000
$ws.Range("M325").Value = 12000
$ws.Range("N325").Value = "`$/saco 25 kilos"
$ws.Range("O325").Value = "Provincia de Cautín"
$ws.Range("P325").Value = 480
$ws.Range("Q325").Value = 25
$ws.Range("R325").Value = "Hortaliza"
